# ---------------------------------------------------------------------------
# Refresh of the "results" text-file query tables (Data > Refresh All) for the
# Maleinsaeure Ragone-plot workbook:
#   - the existing simulation results table (G:J, rows 3-12) was re-pulled
#     from results.csv and now only spans 9 data rows (3-11) instead of 10,
#     with updated values; the old row 12 is cleared.
#   - a brand-new query table ("results_9") was added at A34:D43 pulling the
#     same refreshed CSV snapshot.
#   - defined names results_4 .. results_9 were created to track the history
#     of query-table refreshes (mirrors results_1 .. results_3 already in the
#     workbook), with results_9 pointing at the live new range.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ergebnisse")

# --- refreshed simulation results (G3:J11), 9 rows now instead of 10 -------
$results = @(
    @(3, 0.04, 0.0519124, 0.0000103825, 5000.0),
    @(4, 0.06, 0.0819682, 0.0000163936, 5000.0),
    @(5, 0.08, 0.0874916, 0.0000261846, 3341.33),
    @(6, 0.1, 0.0326871, 0.0000440416, 742.188),
    @(7, 0.12, 0.0202101, 0.0000558036, 362.165),
    @(8, 0.14, 0.0139643, 0.0000637749, 218.962),
    @(9, 0.16, 0.009214, 0.0000674652, 136.574),
    @(10, 0.18, 0.00514523, 0.0000638471, 80.5868),
    @(11, 0.2, 0.00167751, 0.0000449715, 37.3016)
)

foreach ($row in $results) {
    $r = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
    $ws.Cells.Item($r, 10).Value = $row[4]
}

# row 12 no longer holds data after the refresh (table shrank by one row) --
$ws.Range("G12:J12").ClearContents()

# --- new query table "results_9" at A34:D43 ---------------------------------
$ws.Range("A34").Value = "Stromstaerke in A"
$ws.Range("B34").Value = "SAC"
$ws.Range("C34").Value = "ASAR"
$ws.Range("D34").Value = "Adsorptionszeit"

$results9 = @(
    @(35, 0.04, 0.0519124, 0.0000103825, 5000.0),
    @(36, 0.06, 0.0819682, 0.0000163936, 5000.0),
    @(37, 0.08, 0.0874916, 0.0000261846, 3341.33),
    @(38, 0.1, 0.0326871, 0.0000440416, 742.188),
    @(39, 0.12, 0.0202101, 0.0000558036, 362.165),
    @(40, 0.14, 0.0139643, 0.0000637749, 218.962),
    @(41, 0.16, 0.009214, 0.0000674652, 136.574),
    @(42, 0.18, 0.00514523, 0.0000638471, 80.5868),
    @(43, 0.2, 0.00167751, 0.0000449715, 37.3016)
)

foreach ($row in $results9) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# match the scientific-notation display format used by the refreshed table's
# "ASAR" column (same format already applied to column I further up)
$ws.Range("C35:C43").NumberFormat = "0.00E+00"

# --- defined names tracking the query-table refresh history ----------------
$ws.Names.Add('results_4', '=Ergebnisse!#REF!')
$ws.Names.Add('results_5', '=Ergebnisse!#REF!')
$ws.Names.Add('results_6', '=Ergebnisse!#REF!')
$ws.Names.Add('results_7', '=Ergebnisse!#REF!')
$ws.Names.Add('results_8', '=Ergebnisse!#REF!')
$ws.Names.Add('results_9', '=Ergebnisse!$A$34:$E$43')

# --- cursor position left where the refresh left it -------------------------
$ws.Range("B13").Select()

Write-Output "done"
